# Update the workbook metadata (Date, Count) and replace the two
# sample quarantine-station rows with the full 20-row list of US
# quarantine station cities.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B8").Value = "2024-01-05T10:12:51-05:00"   # Date

# "Count" (B21) must remain text "20" (not get auto-converted to a
# number) while keeping its existing cell style. Stage the text value
# in a scratch cell formatted as Text, then paste just the value back
# onto B21 so its own formatting/style is left untouched.
$scratch1 = $ws1.Range("Z1")
$scratch1.NumberFormat = "@"
$scratch1.Value = "20"
$scratch1.Copy()
$ws1.Range("B21").PasteSpecial(-4163)   # xlPasteValues
$scratch1.Clear()

# ---- Concepts sheet --------------------------------------------------
$ws2 = $wb.Worksheets.Item("Concepts")

# City code / display pairs, in row order (row 2 is the first data row).
# Code column (B) and Display column (C) differ only when the code has
# no spaces and the display is the human readable, spaced version.
$cities = @(
  @("Anchorage", "Anchorage"),
  @("Atlanta", "Atlanta"),
  @("Boston", "Boston"),
  @("Chicago", "Chicago"),
  @("Dallas", "Dallas"),
  @("Detroit", "Detroit"),
  @("ElPaso", "El Paso"),
  @("Honolulu", "Honolulu"),
  @("Houston", "Houston"),
  @("LosAngeles", "Los Angeles"),
  @("Miami", "Miami"),
  @("Minneapolis", "Minneapolis"),
  @("NewYork", "New York"),
  @("Newark", "Newark"),
  @("Philadelphia", "Philadelphia"),
  @("SanDiego", "San Diego"),
  @("SanFrancisco", "San Francisco"),
  @("SanJuan", "San Juan"),
  @("Seattle", "Seattle"),
  @("WashingtonDC", "Washington, DC")
)

$template = $ws2.Range("A2:D2")
$startRow = 2

for ($i = 0; $i -lt $cities.Length; $i++) {
  $row = $startRow + $i
  if ($row -gt $startRow) {
    # Copy the existing row 2 as a template so the new row picks up the
    # same styles (and the "1" Level value in column A) exactly.
    $dst = $ws2.Range("A$($row):D$($row)")
    $template.Copy($dst)
  }
  $ws2.Range("B$row").Value = $cities[$i][0]
  $ws2.Range("C$row").Value = $cities[$i][1]
}
